$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.187.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.80%  '

$ws.Range("D3").Value = '''3.663.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.97%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''590.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '

$ws.Range("D6").Value = '''178.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.90%  '

$ws.Range("D7").Value = '''3.666.71'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.67%  '

$ws.Range("E8").Value = '  -5.76%  '

$ws.Range("D9").Value = '''0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").Value = '''0.710'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.52%  '

$ws.Range("E11").Value = '  -9.10%  '

$ws.Range("D12").Value = '''55.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.62%  '

$ws.Range("D13").Value = '''0.0000291'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.39%  '

$ws.Range("D14").Value = '''10.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.56%  '

$ws.Range("D15").Value = '''4.235.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.04%  '

$ws.Range("D16").Value = '''3.652.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.36%  '

$ws.Range("D17").Value = '''19.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.90%  '

$ws.Range("E18").Value = '  -2.22%  '

$ws.Range("E19").Value = '  -7.13%  '

$ws.Range("D20").Value = '''12.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.23%  '

$ws.Range("D21").Value = '''67.941.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.89%  '

$ws.Range("D22").Value = '''408.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.89%  '

$ws.Range("E23").Value = '  -5.02%  '

$ws.Range("D24").Value = '''88.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.48%  '

$ws.Range("D25").Value = '''2.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.98%  '

$ws.Range("D26").Value = '''12.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.88%  '

$ws.Range("D27").Value = '''10.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.73%  '

$ws.Range("D28").Value = '''3.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.78%  '

$ws.Range("D29").Value = '''6.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.48%  '

$ws.Range("D30").Value = '''9.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.17%  '

$ws.Range("D31").Value = '''32.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.87%  '

$ws.Range("D32").Value = '''7.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -15.72%  '

$ws.Range("D33").Value = '''12.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.88%  '

$ws.Range("E34").Value = '  -6.63%  '

$ws.Range("D35").Value = '''64.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.49%  '

$ws.Range("D36").Value = '''601.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.91%  '

$ws.Range("D37").Value = '''42.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.80%  '

$ws.Range("D38").Value = '''0.0₃0883'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.90%  '

$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("D40").Value = '''0.396'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.04%  '

$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("E42").Value = '  -6.68%  '

$ws.Range("D43").Value = '''3.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.87%  '

$ws.Range("D44").Value = '''2.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.34%  '

$ws.Range("E45").Value = '  -7.27%  '

$ws.Range("D46").Value = '''2.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.99%  '

$ws.Range("D47").Value = '''2.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.55%  '

$ws.Range("E48").Value = '  -6.64%  '

$ws.Range("E49").Value = '  -12.02%  '

$ws.Range("D50").Value = '''2.708.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.25%  '

$ws.Range("E51").Value = '  -5.90%  '
